$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 2: "TextBox 5" (Elinor / Nofet / Lidor credits box) ---
# Add a new first paragraph "Moderators: Prof. Kobi Gal, Yakir Ben-Aliz"
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Place the new sentence before the existing first paragraph's text and
# press Enter, splitting it into its own paragraph (mirrors typing at the
# very start of the "Elinor Avraham | ..." line).
[void]$tr.InsertBefore("Moderators: Prof. Kobi Gal, Yakir Ben-Aliz`r")

# Re-anchor on the freshly created first paragraph. It already inherited
# size/bold/Tahoma formatting (with the full panose/pitchFamily/charset
# attributes) from the paragraph it was typed in front of, so only the
# language needs to change to English for the newly typed text.
$full = $sh.TextFrame.TextRange
$para1 = $full.Paragraphs(1, 1)
$para1.LanguageID = "en-US"

# Break the line up into the same run spans as the authored deck:
#   "Moderators: Prof. Kobi Gal, " | "Yakir" | " Ben-" | "Aliz"
# Re-applying Bold (already true) forces each span to materialise as its
# own <a:r> run without disturbing the inherited <a:latin>/<a:ea>/<a:cs>
# typeface attributes (panose/pitchFamily/charset survive intact).
$run1 = $para1.Characters(1, 28)
$run2 = $para1.Characters(29, 5)
$run3 = $para1.Characters(34, 5)
$run4 = $para1.Characters(39, 4)
$run1.Font.Bold = $true
$run2.Font.Bold = $true
$run3.Font.Bold = $true
$run4.Font.Bold = $true

# Resize/reposition the textbox to its final authored geometry (the shape
# uses spAutoFit, so pin down the exact box produced by the edit). Only
# Top/Height actually move; Left/Width are unchanged by the edit.
$sh.Top = 5320246 / 12700.0
$sh.Height = 2031325 / 12700.0

# --- Shape 3: "Picture 6" (decorative divider line) moves up slightly ---
$pic = $s.Shapes.Item(3)
$pic.Top = 4813708 / 12700.0
